$wb = $excel.ActiveWorkbook

# --- Sheet: labor_incmon_imp_stochastic_reg ---
$ws = $wb.Worksheets.Item("labor_incmon_imp_stochastic_reg")
$ws.Range("G3").Value = 2062039.5297215781
$ws.Range("H3").Value = 300000
$ws.Range("I3").Value = 500394.75
$ws.Range("J3").Value = 1400000
$ws.Range("K3").Value = 2654159
$ws.Range("L3").Value = 3862165

# --- Sheet: labor_jubpenimp_stochastic_reg ---
$ws = $wb.Worksheets.Item("labor_jubpenimp_stochastic_reg")
$ws.Range("G3").Value = 1518214.4247522387

# --- Sheet: nonlabor_imp_stochastic_reg ---
$ws = $wb.Worksheets.Item("nonlabor_imp_stochastic_reg")
$ws.Range("G3").Value = 435286.01169894467
$ws.Range("K3").Value = 486977.59375
$ws.Range("L3").Value = 834000

# --- Sheet: labor_beneimp_stochastic_reg ---
$ws = $wb.Worksheets.Item("labor_beneimp_stochastic_reg")
$ws.Range("G3").Value = 898015.88863384188
$ws.Range("I3").Value = 200000
$ws.Range("J3").Value = 328258.90625
$ws.Range("K3").Value = 922541.125
$ws.Range("L3").Value = 2001578.75
